$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart
$ws.Range("H19").Value = 13394.733
$ws.Range("I19").Value = 101.625
$ws.Range("J19").Value = 28586.857
$ws.Range("K19").Value = 101.625
$ws.Range("L19").Value = 28586.857
$ws.Range("M19").Value = 73.375
$ws.Range("N19").Value = -28936.857

# Row 70: Consecrating Congregation
$ws.Range("H70").Value = 1387490.5
$ws.Range("J70").Value = 5004
$ws.Range("L70").Value = 15012
$ws.Range("N70").Value = -15552

# Row 73: Curbing the Contagion (L)
$ws.Range("H73").Value = 1387490.5
$ws.Range("J73").Value = 5004
$ws.Range("L73").Value = 15012
$ws.Range("N73").Value = -16884

# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 1420838.9
$ws.Range("I80").Value = 1748428.9
$ws.Range("K80").Value = 5245286.699999999
$ws.Range("M80").Value = -5244288.699999999

# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 1420838.9
$ws.Range("I83").Value = 1748428.9
$ws.Range("K83").Value = 15735860.1
$ws.Range("M83").Value = -15730868.1

# Row 94: Magic Beans
$ws.Range("H94").Value = 587.5
$ws.Range("I94").Value = 587.5
$ws.Range("K94").Value = 587.5
$ws.Range("M94").Value = -136.5

# Row 100: Asking for a Friend
$ws.Range("H100").Value = 3693.5652
$ws.Range("J100").Value = 975
$ws.Range("L100").Value = 975
$ws.Range("N100").Value = -2057

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 5924.5
$ws.Range("I132").Value = 5639.2354
$ws.Range("J132").Value = 8349.25
$ws.Range("K132").Value = 16917.7062
$ws.Range("L132").Value = 25047.75
$ws.Range("M132").Value = -14387.7062
$ws.Range("N132").Value = -30107.75

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 4176.8335
$ws.Range("I137").Value = 3705.7058
$ws.Range("J137").Value = 4792.923
$ws.Range("K137").Value = 11117.1174
$ws.Range("L137").Value = 14378.769
$ws.Range("M137").Value = -8567.117400000001
$ws.Range("N137").Value = -19478.769

# Row 138: All-night Crafting
$ws.Range("H138").Value = 5950.8057
$ws.Range("I138").Value = 5030
$ws.Range("K138").Value = 15090
$ws.Range("M138").Value = -9950

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 1440.2307
$ws.Range("I2").Value = 410.8889
$ws.Range("K2").Value = 410.8889
$ws.Range("M2").Value = -297.8889

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 31035.867
$ws.Range("I32").Value = 21606
$ws.Range("J32").Value = 45180.668
$ws.Range("K32").Value = 21606
$ws.Range("L32").Value = 45180.668
$ws.Range("M32").Value = -21319
$ws.Range("N32").Value = -45754.668

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 374216
$ws.Range("I61").Value = 3282.4211
$ws.Range("J61").Value = 1255183.2
$ws.Range("K61").Value = 3282.4211
$ws.Range("L61").Value = 1255183.2
$ws.Range("M61").Value = -3070.4211
$ws.Range("N61").Value = -1255607.2

# Row 113: Catching an Earful
$ws.Range("H113").Value = 42000
$ws.Range("J113").Value = 42000
$ws.Range("L113").Value = 42000
$ws.Range("N113").Value = -50678

# Row 116: No Scope
$ws.Range("H116").Value = 1440.2307
$ws.Range("I116").Value = 410.8889
$ws.Range("K116").Value = 410.8889
$ws.Range("M116").Value = 1883.1111

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 447023.97
$ws.Range("I132").Value = 272340.06
$ws.Range("J132").Value = 1254937.1
$ws.Range("K132").Value = 817020.1799999999
$ws.Range("L132").Value = 3764811.3
$ws.Range("M132").Value = -814490.1799999999
$ws.Range("N132").Value = -3769871.3

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 374216
$ws.Range("I136").Value = 3282.4211
$ws.Range("J136").Value = 1255183.2
$ws.Range("K136").Value = 9847.263300000001
$ws.Range("L136").Value = 3765549.6
$ws.Range("M136").Value = -7297.263300000001
$ws.Range("N136").Value = -3770649.6

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 1440.2307
$ws.Range("I3").Value = 410.8889
$ws.Range("K3").Value = 410.8889
$ws.Range("M3").Value = -296.8889

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 4840.2666
$ws.Range("I134").Value = 4537.1816
$ws.Range("K134").Value = 13611.5448
$ws.Range("M134").Value = -11076.5448

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 5107.722
$ws.Range("J31").Value = 5281.5
$ws.Range("L31").Value = 5281.5
$ws.Range("N31").Value = -5871.5

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 5107.722
$ws.Range("J34").Value = 5281.5
$ws.Range("L34").Value = 5281.5
$ws.Range("N34").Value = -5685.5

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 561034.8
$ws.Range("I132").Value = 5759.143
$ws.Range("K132").Value = 17277.429
$ws.Range("M132").Value = -14747.429

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 4911.44
$ws.Range("I134").Value = 4194.3687
$ws.Range("J134").Value = 7182.1665
$ws.Range("K134").Value = 12583.1061
$ws.Range("L134").Value = 21546.4995
$ws.Range("M134").Value = -10048.1061
$ws.Range("N134").Value = -26616.4995

# Row 141: No Greater Treasure
$ws.Range("H141").Value = 255286.14
$ws.Range("J141").Value = 255286.14
$ws.Range("L141").Value = 255286.14
$ws.Range("N141").Value = -265646.14

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 726.375
$ws.Range("J2").Value = 1133
$ws.Range("L2").Value = 6798
$ws.Range("N2").Value = -7024

# Row 5: What a Sap
$ws.Range("H5").Value = 668.5714
$ws.Range("I5").Value = 632.36365
$ws.Range("J5").Value = 801.3333
$ws.Range("K5").Value = 1897.09095
$ws.Range("L5").Value = 2403.9999
$ws.Range("M5").Value = -1785.09095
$ws.Range("N5").Value = -2627.9999

# Row 23: Sweet Smell of Success
$ws.Range("H23").Value = 153
$ws.Range("I23").Value = 156.75
$ws.Range("K23").Value = 470.25
$ws.Range("M23").Value = -235.25

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 5130954
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5130954
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 15392862
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -15397202

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 668.5714
$ws.Range("I135").Value = 632.36365
$ws.Range("J135").Value = 801.3333
$ws.Range("K135").Value = 5691.27285
$ws.Range("L135").Value = 7211.9997
$ws.Range("M135").Value = -3156.27285
$ws.Range("N135").Value = -12281.9997

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 6199.8
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 6199.8
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 18599.4
$ws.Range("M137").Value = $null
$ws.Range("N137").Value = -28799.4

# Row 141: Ocean Explosion
$ws.Range("H141").Value = 14522.25
$ws.Range("I141").Value = 14522.25
$ws.Range("K141").Value = 43566.75
$ws.Range("M141").Value = -38386.75

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 7632.6665
$ws.Range("J70").Value = 7632.6665
$ws.Range("L70").Value = 7632.6665
$ws.Range("N70").Value = -8172.6665

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 7632.6665
$ws.Range("J73").Value = 7632.6665
$ws.Range("L73").Value = 7632.6665
$ws.Range("N73").Value = -9504.666499999999

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 3498.0833
$ws.Range("I80").Value = 2995.8
$ws.Range("J80").Value = 3856.8572
$ws.Range("K80").Value = 2995.8
$ws.Range("L80").Value = 3856.8572
$ws.Range("M80").Value = -1997.8
$ws.Range("N80").Value = -5852.8572

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 3498.0833
$ws.Range("I83").Value = 2995.8
$ws.Range("J83").Value = 3856.8572
$ws.Range("K83").Value = 14979
$ws.Range("L83").Value = 19284.286
$ws.Range("M83").Value = -9987
$ws.Range("N83").Value = -29268.286

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 7414.1943
$ws.Range("I126").Value = 14122
$ws.Range("J126").Value = 5178.2593
$ws.Range("K126").Value = 42366
$ws.Range("L126").Value = 15534.7779
$ws.Range("M126").Value = -39896
$ws.Range("N126").Value = -20474.7779

# Row 132: On Board for Lar
$ws.Range("H132").Value = 406810.72
$ws.Range("I132").Value = 561070.75
$ws.Range("K132").Value = 1683212.25
$ws.Range("M132").Value = -1680682.25

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 435857.47
$ws.Range("J61").Value = 526.6667
$ws.Range("L61").Value = 526.6667
$ws.Range("N61").Value = -930.6667

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 2821.9
$ws.Range("I82").Value = 2707.2
$ws.Range("K82").Value = 2707.2
$ws.Range("M82").Value = -2346.2

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 2821.9
$ws.Range("I85").Value = 2707.2
$ws.Range("K85").Value = 2707.2
$ws.Range("M85").Value = -1459.2

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 3214.5715
$ws.Range("I93").Value = 3214.5715
$ws.Range("K93").Value = 3214.5715
$ws.Range("M93").Value = -1966.5715

# Row 113: Peace in Rest
$ws.Range("H113").Value = 435857.47
$ws.Range("J113").Value = 526.6667
$ws.Range("L113").Value = 526.6667
$ws.Range("N113").Value = -4866.6667

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 169270.89
$ws.Range("I132").Value = 297808.78
$ws.Range("J132").Value = 7408.3335
$ws.Range("K132").Value = 893426.3400000001
$ws.Range("L132").Value = 22225.0005
$ws.Range("M132").Value = -890896.3400000001
$ws.Range("N132").Value = -27285.0005

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 29419682
$ws.Range("I136").Value = 45463256
$ws.Range("J136").Value = 6468.4165
$ws.Range("K136").Value = 136389768
$ws.Range("L136").Value = 19405.2495
$ws.Range("M136").Value = -136387218
$ws.Range("N136").Value = -24505.2495

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 2208.3447
$ws.Range("I122").Value = 1899.3125
$ws.Range("K122").Value = 5697.9375
$ws.Range("M122").Value = -3247.9375

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 507132.6
$ws.Range("I132").Value = 722654.4399999999
$ws.Range("J132").Value = 4248.3335
$ws.Range("K132").Value = 2167963.32
$ws.Range("L132").Value = 12745.0005
$ws.Range("M132").Value = -2165433.32
$ws.Range("N132").Value = -17805.0005

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 11474.685
$ws.Range("I136").Value = 16318.083
$ws.Range("J136").Value = 3171.7144
$ws.Range("K136").Value = 48954.249
$ws.Range("L136").Value = 9515.143199999999
$ws.Range("M136").Value = -46404.249
$ws.Range("N136").Value = -14615.1432
